$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18089.564
$ws.Range("I32").Value = 18013.146
$ws.Range("K32").Value = 18013.146
$ws.Range("M32").Value = -17726.146

$ws.Range("H74").Value = 957.5526
$ws.Range("I74").Value = 957.5526
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 957.5526
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -83.55259999999998
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 957.5526
$ws.Range("I77").Value = 957.5526
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4787.763
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -419.7629999999999
$ws.Range("N77").ClearContents()

$ws.Range("H122").Value = 1900.3448
$ws.Range("I122").Value = 1826.05
$ws.Range("J122").Value = 2065.4443
$ws.Range("K122").Value = 5478.15
$ws.Range("L122").Value = 6196.3329
$ws.Range("M122").Value = -3028.15
$ws.Range("N122").Value = -11096.3329

$ws.Range("H132").Value = 8461.706
$ws.Range("I132").Value = 11280.5
$ws.Range("J132").Value = 3293.9167
$ws.Range("K132").Value = 33841.5
$ws.Range("L132").Value = 9881.750100000001
$ws.Range("M132").Value = -31311.5
$ws.Range("N132").Value = -14941.7501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 358.5263
$ws.Range("I94").Value = 358.5263
$ws.Range("K94").Value = 358.5263
$ws.Range("M94").Value = 92.47370000000001

$ws.Range("H105").Value = 3201.1904
$ws.Range("I105").Value = 2111.4546
$ws.Range("K105").Value = 2111.4546
$ws.Range("M105").Value = -364.4546

$ws.Range("H134").Value = 5974.759
$ws.Range("I134").Value = 7630.9473
$ws.Range("J134").Value = 2828
$ws.Range("K134").Value = 22892.8419
$ws.Range("L134").Value = 8484
$ws.Range("M134").Value = -20357.8419
$ws.Range("N134").Value = -13554

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 669.8148
$ws.Range("I107").Value = 661.0303
$ws.Range("J107").Value = 683.619
$ws.Range("K107").Value = 661.0303
$ws.Range("L107").Value = 683.619
$ws.Range("M107").Value = 1258.9697
$ws.Range("N107").Value = -4523.619

$ws.Range("H132").Value = 4811075
$ws.Range("I132").Value = 3073.7
$ws.Range("K132").Value = 9221.099999999999
$ws.Range("M132").Value = -6691.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1118.2424
$ws.Range("I5").Value = 235.33333
$ws.Range("J5").Value = 2663.3333
$ws.Range("K5").Value = 705.99999
$ws.Range("L5").Value = 7989.999899999999
$ws.Range("M5").Value = -593.99999
$ws.Range("N5").Value = -8213.999899999999

$ws.Range("H41").Value = 680
$ws.Range("J41").Value = 750
$ws.Range("L41").Value = 2250
$ws.Range("N41").Value = -2926

$ws.Range("H62").Value = 3044.9092
$ws.Range("I62").Value = 2499.25
$ws.Range("K62").Value = 7497.75
$ws.Range("M62").Value = -6811.75

$ws.Range("H64").Value = 2148.5
$ws.Range("I64").Value = 1446
$ws.Range("J64").Value = 2499.75
$ws.Range("K64").Value = 4338
$ws.Range("L64").Value = 7499.25
$ws.Range("M64").Value = -4068
$ws.Range("N64").Value = -8039.25

$ws.Range("H65").Value = 3044.9092
$ws.Range("I65").Value = 2499.25
$ws.Range("K65").Value = 22493.25
$ws.Range("M65").Value = -19061.25

$ws.Range("H67").Value = 2148.5
$ws.Range("I67").Value = 1446
$ws.Range("J67").Value = 2499.75
$ws.Range("K67").Value = 4338
$ws.Range("L67").Value = 7499.25
$ws.Range("M67").Value = -3402
$ws.Range("N67").Value = -9371.25

$ws.Range("H131").Value = 631.95
$ws.Range("I131").Value = 259.96875
$ws.Range("J131").Value = 807
$ws.Range("K131").Value = 779.90625
$ws.Range("L131").Value = 2421
$ws.Range("M131").Value = 4260.09375
$ws.Range("N131").Value = -12501

$ws.Range("H132").Value = 1812.4286
$ws.Range("I132").Value = 829
$ws.Range("J132").Value = 2550
$ws.Range("K132").Value = 7461
$ws.Range("L132").Value = 22950
$ws.Range("M132").Value = -4931
$ws.Range("N132").Value = -28010

$ws.Range("H135").Value = 1118.2424
$ws.Range("I135").Value = 235.33333
$ws.Range("J135").Value = 2663.3333
$ws.Range("K135").Value = 2117.99997
$ws.Range("L135").Value = 23969.9997
$ws.Range("M135").Value = 417.0000300000002
$ws.Range("N135").Value = -29039.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2162.6316
$ws.Range("I97").Value = 2020.7142
$ws.Range("J97").Value = 2560
$ws.Range("K97").Value = 2020.7142
$ws.Range("L97").Value = 2560
$ws.Range("M97").Value = -1524.7142
$ws.Range("N97").Value = -3552

$ws.Range("H102").Value = 1542.55
$ws.Range("I102").Value = 1184.3846
$ws.Range("J102").Value = 2207.7144
$ws.Range("K102").Value = 1184.3846
$ws.Range("L102").Value = 2207.7144
$ws.Range("M102").Value = 437.6153999999999
$ws.Range("N102").Value = -5451.7144

$ws.Range("H131").Value = 22000
$ws.Range("J131").Value = 22000
$ws.Range("L131").Value = 22000
$ws.Range("N131").Value = -32080

$ws.Range("H132").Value = 4562.6113
$ws.Range("I132").Value = 5066.2856
$ws.Range("K132").Value = 15198.8568
$ws.Range("M132").Value = -12668.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 11495263
$ws.Range("I61").Value = 916.9091
$ws.Range("J61").Value = 47620350
$ws.Range("K61").Value = 916.9091
$ws.Range("L61").Value = 47620350
$ws.Range("M61").Value = -714.9091
$ws.Range("N61").Value = -47620754

$ws.Range("H113").Value = 11495263
$ws.Range("I113").Value = 916.9091
$ws.Range("J113").Value = 47620350
$ws.Range("K113").Value = 916.9091
$ws.Range("L113").Value = 47620350
$ws.Range("M113").Value = 1253.0909
$ws.Range("N113").Value = -47624690

$ws.Range("H132").Value = 9116
$ws.Range("I132").Value = 15926.571
$ws.Range("K132").Value = 47779.713
$ws.Range("M132").Value = -45249.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 24666.666
$ws.Range("J112").Value = 24666.666
$ws.Range("L112").Value = 24666.666
$ws.Range("N112").Value = -27620.666

$ws.Range("H122").Value = 1708.4314
$ws.Range("I122").Value = 1627
$ws.Range("J122").Value = 1946.4615
$ws.Range("K122").Value = 4881
$ws.Range("L122").Value = 5839.3845
$ws.Range("M122").Value = -2431
$ws.Range("N122").Value = -10739.3845

$ws.Range("H126").Value = 27783096
$ws.Range("I126").Value = 43484708
$ws.Range("K126").Value = 130454124
$ws.Range("M126").Value = -130451654

$ws.Range("H132").Value = 1741.5416
$ws.Range("I132").Value = 1907.5358
$ws.Range("J132").Value = 1509.15
$ws.Range("K132").Value = 5722.607400000001
$ws.Range("L132").Value = 4527.450000000001
$ws.Range("M132").Value = -3192.607400000001
$ws.Range("N132").Value = -9587.450000000001

$ws.Range("H133").Value = 42225.75
$ws.Range("J133").Value = 42225.75
$ws.Range("L133").Value = 42225.75
$ws.Range("N133").Value = -52345.75
